$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New section "Crud Usuario": header row 17 (merged, bold/centered like
# the existing "Login" / "Recupera contraseña" section headers) followed
# by the column-title row 18 and a new table Tabla24 over A18:D22.
# ---------------------------------------------------------------------
[void]$ws.Range("A17:D17").Merge()
$ws.Range("A17").Value = "Crud Usuario"
$ws.Range("A10:D10").Copy()
[void]$ws.Range("A17:D17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A18").Value = "Valores De Prueba"
$ws.Range("B18").Value = "Accion"
$ws.Range("C18").Value = "Resultado Esperado"
$ws.Range("D18").Value = "Resultado Obtenido"

$lo1 = $ws.ListObjects.Add(1, $ws.Range("A18:D22"), 0, 1)
$lo1.Name = "Tabla24"
$lo1.TableStyle = "TableStyleMedium1"

# ---------------------------------------------------------------------
# New section "Crud Productos": header row 24 (merged, bold/centered)
# followed by the column-title row 25 and a new table Tabla245 over
# A25:D29.
# ---------------------------------------------------------------------
[void]$ws.Range("A24:D24").Merge()
$ws.Range("A24").Value = "Crud Productos"
$ws.Range("A10:D10").Copy()
[void]$ws.Range("A24:D24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A25").Value = "Valores De Prueba"
$ws.Range("B25").Value = "Accion"
$ws.Range("C25").Value = "Resultado Esperado"
$ws.Range("D25").Value = "Resultado Obtenido"

$lo2 = $ws.ListObjects.Add(1, $ws.Range("A25:D29"), 0, 1)
$lo2.Name = "Tabla245"
$lo2.TableStyle = "TableStyleMedium1"

# ---------------------------------------------------------------------
# Selection ends up on B33, matching the saved workbook view state.
# ---------------------------------------------------------------------
[void]$ws.Range("B33").Select()
